# Update the "dSF" column (column F) values for the castillo_luis sheet.
# These are the updated delta-stock-final values coming from a repull of
# the source data (see commit message: "repull data, push all data, mean
# calculation"). Only column F changes; all other columns/values are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -2
    4  = 3
    6  = -3
    7  = -2
    8  = -1
    9  = 3
    10 = -2
    11 = -5
    12 = 4
    13 = 1
    14 = 2
    15 = 2
    17 = 5
    18 = 2
    19 = 3
    20 = -3
    22 = 1
    23 = 3
    25 = -4
    26 = -2
    27 = -7
    28 = -5
    30 = -4
    31 = -2
    32 = 1
    33 = 8
    34 = -6
    36 = 2
    37 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
